$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

$ws.Range("A3").Value = "jacksonions0502"
$ws.Range("B3").Value = "jackisions0502"
$ws.Range("A4").Value = "averyions0502"
$ws.Range("B4").Value = "wyattions0502"
